$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.818.20'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '2.498.35'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'" + '322.72'
$ws.Range('D6').Value = "'" + '108.99'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('D7').Value = "'" + '0.524'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +1.86%  '
$ws.Range('D10').Value = "'" + '40.24'
$ws.Range('E10').Value = '  +5.53%  '
$ws.Range('D11').Value = "'" + '0.0814'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'" + '0.124'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = "'" + '18.94'
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = '2.888.94'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '2.497.34'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').Value = "'" + '0.848'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').Value = '47.693.58'
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('D19').Value = "'" + '13.15'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('E21').Value = '  +12.21%  '
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = "'" + '70.76'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = "'" + '247.81'
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('D25').Value = "'" + '2.56'
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = "'" + '25.89'
$ws.Range('E27').Value = '  -1.17%  '
$ws.Range('D28').Value = "'" + '9.97'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'" + '2.19'
$ws.Range('E29').Value = '  -4.37%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = "'" + '0.139'
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = "'" + '35.03'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').Value = "'" + '49.82'
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('D33').Value = "'" + '19.84'
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('D34').Value = "'" + '5.36'
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  -1.33%  '
$ws.Range('D38').Value = "'" + '4.66'
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').Value = "'" + '0.112'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').Value = "'" + '22.29'
$ws.Range('E41').Value = '  +5.26%  '
$ws.Range('E42').Value = '  -1.47%  '
$ws.Range('D43').Value = "'" + '119.33'
$ws.Range('E43').Value = '  -2.62%  '
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Value = '2.004.33'
$ws.Range('E45').Value = '  +1.78%  '
$ws.Range('D46').Value = "'" + '3.05'
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('E47').Value = '  -3.50%  '
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  -2.59%  '
$ws.Range('E51').Value = '  +2.81%  '
